$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10810.75
$ws.Range("I51").Value = 8994
$ws.Range("J51").Value = 11416.333
$ws.Range("K51").Value = 8994
$ws.Range("L51").Value = 11416.333
$ws.Range("M51").Value = -8510
$ws.Range("N51").Value = -12384.333
$ws.Range("H64").Value = 4674.3687
$ws.Range("I64").Value = 3453.25
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 3453.25
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -3205.25
$ws.Range("N64").Value = -5496
$ws.Range("H67").Value = 4674.3687
$ws.Range("I67").Value = 3453.25
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 3453.25
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -2595.25
$ws.Range("N67").Value = -6716
$ws.Range("H100").Value = 982.8
$ws.Range("I100").Value = 965.625
$ws.Range("K100").Value = 965.625
$ws.Range("M100").Value = -424.625
$ws.Range("H113").Value = 5730.778
$ws.Range("I113").Value = 4888.5
$ws.Range("K113").Value = 4888.5
$ws.Range("M113").Value = -1634.5
$ws.Range("H137").Value = 778212.5600000001
$ws.Range("I137").Value = 1253511.8
$ws.Range("K137").Value = 3760535.4
$ws.Range("M137").Value = -3757985.4
$ws.Range("H138").Value = 2948.7937
$ws.Range("I138").Value = 1699.4
$ws.Range("J138").Value = 3016.6956
$ws.Range("K138").Value = 5098.200000000001
$ws.Range("L138").Value = 9050.086800000001
$ws.Range("M138").Value = 41.79999999999927
$ws.Range("N138").Value = -19330.0868

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5231447.5
$ws.Range("I32").Value = 7262889
$ws.Range("J32").Value = 39986
$ws.Range("K32").Value = 7262889
$ws.Range("L32").Value = 39986
$ws.Range("M32").Value = -7262602
$ws.Range("N32").Value = -40560
$ws.Range("H74").Value = 11374774
$ws.Range("I74").Value = 19235824
$ws.Range("K74").Value = 19235824
$ws.Range("M74").Value = -19234950
$ws.Range("H77").Value = 11374774
$ws.Range("I77").Value = 19235824
$ws.Range("K77").Value = 96179120
$ws.Range("M77").Value = -96174752
$ws.Range("H102").Value = 14651.091
$ws.Range("I102").Value = 16057.5
$ws.Range("J102").Value = 587
$ws.Range("K102").Value = 16057.5
$ws.Range("L102").Value = 587
$ws.Range("M102").Value = -14435.5
$ws.Range("N102").Value = -3831
$ws.Range("H103").Value = 72617.336
$ws.Range("J103").Value = 72617.336
$ws.Range("L103").Value = 72617.336
$ws.Range("N103").Value = -74961.336
$ws.Range("H110").Value = 1318.3334
$ws.Range("I110").Value = 1265.4546
$ws.Range("K110").Value = 1265.4546
$ws.Range("M110").Value = 779.5454
$ws.Range("H122").Value = 3400.5
$ws.Range("I122").Value = 2440.8
$ws.Range("K122").Value = 7322.400000000001
$ws.Range("M122").Value = -4872.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 12762.6
$ws.Range("I26").Value = 12762.6
$ws.Range("K26").Value = 12762.6
$ws.Range("M26").Value = -12470.6
$ws.Range("H96").Value = 19330.46
$ws.Range("I96").Value = 11858
$ws.Range("J96").Value = 109000
$ws.Range("K96").Value = 11858
$ws.Range("L96").Value = 109000
$ws.Range("M96").Value = -9112
$ws.Range("N96").Value = -114492
$ws.Range("H105").Value = 3487.5
$ws.Range("I105").Value = 1900
$ws.Range("J105").Value = 4016.6667
$ws.Range("K105").Value = 1900
$ws.Range("L105").Value = 4016.6667
$ws.Range("M105").Value = -153
$ws.Range("N105").Value = -7510.6667
$ws.Range("H131").Value = 99999.5
$ws.Range("I131").Value = 99999
$ws.Range("K131").Value = 99999
$ws.Range("M131").Value = -94959
$ws.Range("H134").Value = 3325253.2
$ws.Range("I134").Value = 1975.1515
$ws.Range("K134").Value = 5925.4545
$ws.Range("M134").Value = -3390.4545

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 941.8
$ws.Range("I16").Value = 770
$ws.Range("J16").Value = 1199.5
$ws.Range("K16").Value = 770
$ws.Range("L16").Value = 1199.5
$ws.Range("M16").Value = -483
$ws.Range("N16").Value = -1773.5
$ws.Range("H31").Value = 1805508.9
$ws.Range("J31").Value = 5827012
$ws.Range("L31").Value = 5827012
$ws.Range("N31").Value = -5827602
$ws.Range("H34").Value = 1805508.9
$ws.Range("J34").Value = 5827012
$ws.Range("L34").Value = 5827012
$ws.Range("N34").Value = -5827416
$ws.Range("H68").Value = 73236.2
$ws.Range("J68").Value = 73236.2
$ws.Range("L68").Value = 73236.2
$ws.Range("N68").Value = -74734.2
$ws.Range("H71").Value = 73236.2
$ws.Range("J71").Value = 73236.2
$ws.Range("L71").Value = 219708.6
$ws.Range("N71").Value = -227196.6
$ws.Range("H113").Value = 941.8
$ws.Range("I113").Value = 770
$ws.Range("J113").Value = 1199.5
$ws.Range("K113").Value = 770
$ws.Range("L113").Value = 1199.5
$ws.Range("M113").Value = 1400
$ws.Range("N113").Value = -5539.5
$ws.Range("H124").Value = 65666
$ws.Range("J124").Value = 65666
$ws.Range("L124").Value = 65666
$ws.Range("N124").Value = -70576

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 822.0909
$ws.Range("I23").Value = 1548.75
$ws.Range("J23").Value = 406.85715
$ws.Range("K23").Value = 4646.25
$ws.Range("L23").Value = 1220.57145
$ws.Range("M23").Value = -4411.25
$ws.Range("N23").Value = -1690.57145
$ws.Range("H132").Value = 1856.8
$ws.Range("I132").Value = 2332.6667
$ws.Range("J132").Value = 1772.8235
$ws.Range("K132").Value = 20994.0003
$ws.Range("L132").Value = 15955.4115
$ws.Range("M132").Value = -18464.0003
$ws.Range("N132").Value = -21015.4115
$ws.Range("H141").Value = 608996
$ws.Range("I141").Value = 757495
$ws.Range("K141").Value = 2272485
$ws.Range("M141").Value = -2267305

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1991.5
$ws.Range("I22").Value = 1987.75
$ws.Range("J22").Value = 1999
$ws.Range("K22").Value = 1987.75
$ws.Range("L22").Value = 1999
$ws.Range("M22").Value = -1692.75
$ws.Range("N22").Value = -2589
$ws.Range("H27").Value = 1991.5
$ws.Range("I27").Value = 1987.75
$ws.Range("J27").Value = 1999
$ws.Range("K27").Value = 1987.75
$ws.Range("L27").Value = 1999
$ws.Range("M27").Value = -1880.75
$ws.Range("N27").Value = -2213
$ws.Range("H55").Value = 71429140
$ws.Range("I55").Value = 76923640
$ws.Range("J55").Value = 499
$ws.Range("K55").Value = 76923640
$ws.Range("L55").Value = 499
$ws.Range("M55").Value = -76923467
$ws.Range("N55").Value = -845
$ws.Range("H82").Value = 2062.8
$ws.Range("I82").Value = 1520.7778
$ws.Range("J82").Value = 2506.2727
$ws.Range("K82").Value = 1520.7778
$ws.Range("L82").Value = 2506.2727
$ws.Range("M82").Value = -1159.7778
$ws.Range("N82").Value = -3228.2727
$ws.Range("H85").Value = 2062.8
$ws.Range("I85").Value = 1520.7778
$ws.Range("J85").Value = 2506.2727
$ws.Range("K85").Value = 1520.7778
$ws.Range("L85").Value = 2506.2727
$ws.Range("M85").Value = -272.7778000000001
$ws.Range("N85").Value = -5002.2727
$ws.Range("H100").Value = 1863.1428
$ws.Range("I100").Value = 1703.6923
$ws.Range("K100").Value = 1703.6923
$ws.Range("M100").Value = -1162.6923

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 17858476
$ws.Range("I107").Value = 26317152
$ws.Range("J107").Value = 1272.6666
$ws.Range("K107").Value = 78951456
$ws.Range("L107").Value = 3817.9998
$ws.Range("M107").Value = -78949536
$ws.Range("N107").Value = -7657.9998
$ws.Range("H132").Value = 3338494.2
$ws.Range("I132").Value = 4360.8
$ws.Range("K132").Value = 13082.4
$ws.Range("M132").Value = -10552.4
$ws.Range("H136").Value = 4134.303
$ws.Range("I136").Value = 3404.4194
$ws.Range("K136").Value = 10213.2582
$ws.Range("M136").Value = -7663.2582
